$d = $word.ActiveDocument
$find = $d.Content.Find

# --- Programa resumido ---
$old1 = "1. Fundamentos da Gestão de Produção2. Visão estratégica da Produção.3. Projeto em Gestão da Produção.4. Planejamento e Controle da Produção"
$new1 = "1. Fundamentos da Gestão de Produção^l2. Visão estratégica da Produção.^l3. Projeto em Gestão da Produção.^l4. Planejamento e Controle da Produção"
$find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)

# --- Programa ---
$old2 = "1 - Fundamentos da gestão de produção: modelo de transformação: inputs, processo de transformação e outputs. Tipos de Processo de Produção2 - Visão estratégica de produção: Papel da função produção. Objetivos de Desempenho. Estratégias de Produção. Ciclo de Vida Produto/Serviço.3 – Projeto em Gestão da Produção: Tipos de Processos. Projeto de Produtos e Serviços. Projeto de Rede de Operações Produtivas. Arranjo Físico.4 - Planejamento e Controle da Produção: Material Requirement Planning (MRP), Manufacturing Resources Planning (MPRII), Enterprise Planning (ERP). Produção Enxuta. Kanban. Just in Time."
$new2 = "1 - Fundamentos da gestão de produção: modelo de transformação: inputs, processo de transformação e outputs. Tipos de Processo de Produção^l2 - Visão estratégica de produção: Papel da função produção. Objetivos de Desempenho. Estratégias de Produção. Ciclo de Vida Produto/Serviço.^l3 – Projeto em Gestão da Produção: Tipos de Processos. Projeto de Produtos e Serviços. Projeto de Rede de Operações Produtivas. Arranjo Físico.^l4 - Planejamento e Controle da Produção: Material Requirement Planning (MRP), Manufacturing Resources Planning (MPRII), Enterprise Planning (ERP). Produção Enxuta. Kanban. Just in Time."
$find.Execute($old2, $true, $false, $false, $false, $false, $true, 1, $false, $new2, 2)

# --- Bibliografia ---
$old3 = "SLACK, N. et al. Administração da Produção. 3 ed. São Paulo: Atlas, 2009.CHASE, R. B. E JACOBS, F.R. Administração da Produção e de Operações. 1 ed. Porto Alegre. Bookman. 2009.CORREA, H.L.; CORREA, C.A. Administração da Produção e Operações. 2 ed. São Paulo. Atlas. 2006"
$new3 = "SLACK, N. et al. Administração da Produção. 3 ed. São Paulo: Atlas, 2009.^lCHASE, R. B. E JACOBS, F.R. Administração da Produção e de Operações. 1 ed. Porto Alegre. Bookman. 2009.^lCORREA, H.L.; CORREA, C.A. Administração da Produção e Operações. 2 ed. São Paulo. Atlas. 2006"
$find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)

Write-Host "Done"
